# Auto-generated PowerShell Excel COM-interop script
# Updates maa:// stat strings and the "update date" cell on Sheet1
# to reflect the latest CI auto-update run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'maa://24702 (94.24), maa://25390 (95.9), maa://36681 (87.01)'
$ws.Range("L2").Value = '*maa://24633 (55.48), *maa://30515 (69.0), *maa://34787 (72.06), ***maa://20792 (11.93), maa://39402 (87.18), ***maa://29083 (27.78)'
$ws.Range("T2").Value = 'maa://22742 (91.61), *maa://20791 (62.86)'
$ws.Range("P3").Value = 'maa://21249 (94.59), maa://26254 (95.83)'
$ws.Range("T3").Value = 'maa://24617 (88.89), **maa://20790 (44.78), ***maa://37170 (18.52)'
$ws.Range("X3").Value = 'maa://27396 (84.42), maa://27484 (96.08), maa://27480 (82.86)'
$ws.Range("T4").Value = 'maa://32509 (97.98), maa://27295 (83.61), maa://22754 (91.67), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range("X4").Value = '**maa://32495 (47.89), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (89.66)'
$ws.Range("AF4").Value = '*maa://30062 (63.04), ***maa://26209 (13.04), *maa://39394 (68.18)'
$ws.Range("H6").Value = 'maa://24370 (96.55)'
$ws.Range("X7").Value = 'maa://22399 (95.17), *maa://22758 (73.33)'
$ws.Range("A8").Value = '更新日期：2024.12.20 08:32:28'
$ws.Range("T10").Value = 'maa://27395 (95.98), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range("X10").Value = 'maa://22301 (97.6), maa://22726 (100.0)'
$ws.Range("T11").Value = 'maa://22747 (93.42), maa://22501 (98.59)'
$ws.Range("P13").Value = 'maa://22676 (91.82), *maa://22583 (75.81), *maa://22500 (57.78)'
$ws.Range("D15").Value = '*maa://22743 (76.96), maa://22734 (83.9), *maa://30808 (64.52), **maa://36048 (33.33)'
$ws.Range("H17").Value = 'maa://22430 (88.65), maa://39599 (86.49)'
$ws.Range("AF19").Value = '*maa://21663 (62.12)'
$ws.Range("X21").Value = 'maa://20110 (86.76), maa://34946 (92.68)'
$ws.Range("AB21").Value = 'maa://21443 (80.11), ***maa://23820 (29.82)'
$ws.Range("L23").Value = 'maa://39756 (94.0), maa://39875 (93.75)'
$ws.Range("X24").Value = 'maa://29988 (86.67), maa://23504 (93.28), **maa://22892 (39.58), *maa://25141 (77.6), maa://36663 (80.6), ***maa://22815 (23.08)'
$ws.Range("H25").Value = '*maa://29063 (73.47), *maa://25311 (74.0), ***maa://22725 (4.84)'
$ws.Range("AF25").Value = 'maa://20108 (96.24), maa://24621 (96.58), maa://36676 (96.67), maa://22771 (85.71), maa://37772 (100.0)'
$ws.Range("X28").Value = 'maa://39929 (89.68), ***maa://39723 (14.29), maa://41749 (90.74)'
$ws.Range("AF29").Value = '*maa://24080 (69.05), ***maa://34960 (8.33), *maa://42865 (78.12)'
$ws.Range("H32").Value = 'maa://21895 (97.27), maa://36667 (98.33), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range("P37").Value = 'maa://21280 (89.05), *maa://21239 (66.67)'
$ws.Range("H39").Value = 'maa://25199 (84.82), maa://36670 (87.06), maa://30434 (89.23), ***maa://25036 (16.0), **maa://44165 (50.0)'
$ws.Range("H45").Value = 'maa://21229 (84.78), maa://30807 (95.45), *maa://22767 (55.0), ***maa://20796 (13.79), *maa://42459 (75.0)'
$ws.Range("H46").Value = 'maa://35931 (92.5), maa://43901 (87.5)'
$ws.Range("H53").Value = 'maa://32534 (93.52), **maa://32434 (34.78)'
